# Auto-generated: apply updated market-price / profit figures per the commit diff.
# Source: scheduled runner refresh of currentAveragePrice* / LevePrice* / LeveProfit* columns
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 1332.3334  # H5: 1016 -> 1332.3334
$ws.Cells.Item(5, 9).Value = 998.5  # I5: 524 -> 998.5
$ws.Cells.Item(5, 11).Value = 998.5  # K5: 524 -> 998.5
$ws.Cells.Item(5, 13).Value = -883.5  # M5: -409 -> -883.5
$ws.Cells.Item(6, 8).Value = 2079.6667  # H6: 1898 -> 2079.6667
$ws.Cells.Item(6, 9).Value = 369.5  # I6: 697.3333 -> 369.5
$ws.Cells.Item(6, 11).Value = 1108.5  # K6: 2091.9999 -> 1108.5
$ws.Cells.Item(6, 13).Value = -996.5  # M6: -1979.9999 -> -996.5
$ws.Cells.Item(12, 8).Value = 1069.6  # H12: 1199.625 -> 1069.6
$ws.Cells.Item(12, 9).Value = 832.5  # I12: 979.2 -> 832.5
$ws.Cells.Item(12, 10).Value = 1425.25  # J12: 1567 -> 1425.25
$ws.Cells.Item(12, 11).Value = 832.5  # K12: 979.2 -> 832.5
$ws.Cells.Item(12, 12).Value = 1425.25  # L12: 1567 -> 1425.25
$ws.Cells.Item(12, 13).Value = -662.5  # M12: -809.2 -> -662.5
$ws.Cells.Item(12, 14).Value = -1765.25  # N12: -1907 -> -1765.25
$ws.Cells.Item(43, 8).Value = 1770.3077  # H43: 1792.9166 -> 1770.3077
$ws.Cells.Item(43, 9).Value = 1744.5  # I43: 1793.6 -> 1744.5
$ws.Cells.Item(43, 11).Value = 1744.5  # K43: 1793.6 -> 1744.5
$ws.Cells.Item(43, 13).Value = -1675.5  # M43: -1724.6 -> -1675.5
$ws.Cells.Item(53, 8).Value = 265.6111  # H53: 508.55554 -> 265.6111
$ws.Cells.Item(53, 9).Value = 327  # I53: 764.3 -> 327
$ws.Cells.Item(53, 11).Value = 327  # K53: 764.3 -> 327
$ws.Cells.Item(53, 13).Value = 310  # M53: -127.3 -> 310
$ws.Cells.Item(70, 8).Value = 1349.4839  # H70: 1361.5 -> 1349.4839
$ws.Cells.Item(70, 9).Value = 2329  # I70: 2999 -> 2329
$ws.Cells.Item(70, 11).Value = 6987  # K70: 8997 -> 6987
$ws.Cells.Item(70, 13).Value = -6717  # M70: -8727 -> -6717
$ws.Cells.Item(73, 8).Value = 1349.4839  # H73: 1361.5 -> 1349.4839
$ws.Cells.Item(73, 9).Value = 2329  # I73: 2999 -> 2329
$ws.Cells.Item(73, 11).Value = 6987  # K73: 8997 -> 6987
$ws.Cells.Item(73, 13).Value = -6051  # M73: -8061 -> -6051
$ws.Cells.Item(86, 8).Value = 3387.8667  # H86: 20153 -> 3387.8667
$ws.Cells.Item(86, 9).Value = 3451.2856  # I86: 3833.7 -> 3451.2856
$ws.Cells.Item(86, 10).Value = 2500  # J86: 101749.5 -> 2500
$ws.Cells.Item(86, 11).Value = 3451.2856  # K86: 3833.7 -> 3451.2856
$ws.Cells.Item(86, 12).Value = 2500  # L86: 101749.5 -> 2500
$ws.Cells.Item(86, 13).Value = -2328.2856  # M86: -2710.7 -> -2328.2856
$ws.Cells.Item(86, 14).Value = -4746  # N86: -103995.5 -> -4746
$ws.Cells.Item(89, 8).Value = 3387.8667  # H89: 20153 -> 3387.8667
$ws.Cells.Item(89, 9).Value = 3451.2856  # I89: 3833.7 -> 3451.2856
$ws.Cells.Item(89, 10).Value = 2500  # J89: 101749.5 -> 2500
$ws.Cells.Item(89, 11).Value = 17256.428  # K89: 19168.5 -> 17256.428
$ws.Cells.Item(89, 12).Value = 12500  # L89: 508747.5 -> 12500
$ws.Cells.Item(89, 13).Value = -11640.428  # M89: -13552.5 -> -11640.428
$ws.Cells.Item(89, 14).Value = -23732  # N89: -519979.5 -> -23732
$ws.Cells.Item(100, 9).Value = 829.63635  # I100: 862.7 -> 829.63635
$ws.Cells.Item(100, 10).Value = 1736.125  # J100: 1663.3529 -> 1736.125
$ws.Cells.Item(100, 11).Value = 829.63635  # K100: 862.7 -> 829.63635
$ws.Cells.Item(100, 12).Value = 1736.125  # L100: 1663.3529 -> 1736.125
$ws.Cells.Item(100, 13).Value = -288.63635  # M100: -321.7 -> -288.63635
$ws.Cells.Item(100, 14).Value = -2818.125  # N100: -2745.3529 -> -2818.125
$ws.Cells.Item(101, 8).Value = 600  # H101: 466 -> 600
$ws.Cells.Item(101, 9).Value = 200  # I101: 199 -> 200
$ws.Cells.Item(101, 11).Value = 600  # K101: 597 -> 600
$ws.Cells.Item(101, 13).Value = 1022  # M101: 1025 -> 1022
$ws.Cells.Item(138, 8).Value = 1365.5366  # H138: 1327.878 -> 1365.5366
$ws.Cells.Item(138, 9).Value = 1250.9744  # I138: 1271.075 -> 1250.9744
$ws.Cells.Item(138, 10).Value = 3599.5  # J138: 3600 -> 3599.5
$ws.Cells.Item(138, 11).Value = 3752.9232  # K138: 3813.225 -> 3752.9232
$ws.Cells.Item(138, 12).Value = 10798.5  # L138: 10800 -> 10798.5
$ws.Cells.Item(138, 13).Value = 1387.0768  # M138: 1326.775 -> 1387.0768
$ws.Cells.Item(138, 14).Value = -21078.5  # N138: -21080 -> -21078.5
$ws.Cells.Item(141, 8).Value = 6066.515  # H141: 6430.9355 -> 6066.515
$ws.Cells.Item(141, 9).Value = 6879.3335  # I141: 7772.7617 -> 6879.3335
$ws.Cells.Item(141, 10).Value = 3899  # J141: 3613.1 -> 3899
$ws.Cells.Item(141, 11).Value = 20638.0005  # K141: 23318.2851 -> 20638.0005
$ws.Cells.Item(141, 12).Value = 11697  # L141: 10839.3 -> 11697
$ws.Cells.Item(141, 13).Value = -15458.0005  # M141: -18138.2851 -> -15458.0005
$ws.Cells.Item(141, 14).Value = -22057  # N141: -21199.3 -> -22057

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 406.6  # H5: 273.25 -> 406.6
$ws.Cells.Item(5, 9).Value = 308.25  # I5: 198 -> 308.25
$ws.Cells.Item(5, 11).Value = 308.25  # K5: 198 -> 308.25
$ws.Cells.Item(5, 13).Value = -196.25  # M5: -86 -> -196.25
$ws.Cells.Item(61, 8).Value = 4387.4375  # H61: 4445.933 -> 4387.4375
$ws.Cells.Item(61, 9).Value = 3315.4614  # I61: 3314.6924 -> 3315.4614
$ws.Cells.Item(61, 10).Value = 9032.666999999999  # J61: 11799 -> 9032.666999999999
$ws.Cells.Item(61, 11).Value = 3315.4614  # K61: 3314.6924 -> 3315.4614
$ws.Cells.Item(61, 12).Value = 9032.666999999999  # L61: 11799 -> 9032.666999999999
$ws.Cells.Item(61, 13).Value = -3103.4614  # M61: -3102.6924 -> -3103.4614
$ws.Cells.Item(61, 14).Value = -9456.666999999999  # N61: -12223 -> -9456.666999999999
$ws.Cells.Item(132, 8).Value = 2082.4358  # H132: 2012.0488 -> 2082.4358
$ws.Cells.Item(132, 9).Value = 2082.4358  # I132: 2012.0488 -> 2082.4358
$ws.Cells.Item(132, 11).Value = 6247.307400000001  # K132: 6036.1464 -> 6247.307400000001
$ws.Cells.Item(132, 13).Value = -3717.307400000001  # M132: -3506.1464 -> -3717.307400000001
$ws.Cells.Item(136, 8).Value = 4387.4375  # H136: 4445.933 -> 4387.4375
$ws.Cells.Item(136, 9).Value = 3315.4614  # I136: 3314.6924 -> 3315.4614
$ws.Cells.Item(136, 10).Value = 9032.666999999999  # J136: 11799 -> 9032.666999999999
$ws.Cells.Item(136, 11).Value = 9946.3842  # K136: 9944.0772 -> 9946.3842
$ws.Cells.Item(136, 12).Value = 27098.001  # L136: 35397 -> 27098.001
$ws.Cells.Item(136, 13).Value = -7396.3842  # M136: -7394.0772 -> -7396.3842
$ws.Cells.Item(136, 14).Value = -32198.001  # N136: -40497 -> -32198.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 406.6  # H4: 273.25 -> 406.6
$ws.Cells.Item(4, 9).Value = 308.25  # I4: 198 -> 308.25
$ws.Cells.Item(4, 11).Value = 308.25  # K4: 198 -> 308.25
$ws.Cells.Item(4, 13).Value = -193.25  # M4: -83 -> -193.25
$ws.Cells.Item(47, 8).Value = 170339.5  # H47: 89999 -> 170339.5
$ws.Cells.Item(47, 10).Value = 170339.5  # J47: 89999 -> 170339.5
$ws.Cells.Item(47, 12).Value = 170339.5  # L47: 89999 -> 170339.5
$ws.Cells.Item(47, 14).Value = -171379.5  # N47: -91039 -> -171379.5
$ws.Cells.Item(109, 8).Value = 59999  # H109: 150000 -> 59999
$ws.Cells.Item(109, 10).Value = 59999  # J109: 150000 -> 59999
$ws.Cells.Item(109, 12).Value = 59999  # L109: 150000 -> 59999
$ws.Cells.Item(109, 14).Value = -62773  # N109: -152774 -> -62773
$ws.Cells.Item(134, 8).Value = 5983.4814  # H134: 5314.4194 -> 5983.4814
$ws.Cells.Item(134, 9).Value = 5829.0386  # I134: 5158.2666 -> 5829.0386
$ws.Cells.Item(134, 11).Value = 17487.1158  # K134: 15474.7998 -> 17487.1158
$ws.Cells.Item(134, 13).Value = -14952.1158  # M134: -12939.7998 -> -14952.1158

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 61.375  # H7: 58.6875 -> 61.375
$ws.Cells.Item(7, 9).Value = 29  # I7: 23.625 -> 29
$ws.Cells.Item(7, 11).Value = 29  # K7: 23.625 -> 29
$ws.Cells.Item(7, 13).Value = 84  # M7: 89.375 -> 84
$ws.Cells.Item(22, 8).Value = 1194.6666  # H22: 1303.6296 -> 1194.6666
$ws.Cells.Item(22, 9).Value = 1207.4  # I22: 1363 -> 1207.4
$ws.Cells.Item(22, 10).Value = 1181.9333  # J22: 1248.5 -> 1181.9333
$ws.Cells.Item(22, 11).Value = 1207.4  # K22: 1363 -> 1207.4
$ws.Cells.Item(22, 12).Value = 1181.9333  # L22: 1248.5 -> 1181.9333
$ws.Cells.Item(22, 13).Value = -857.4000000000001  # M22: -1013 -> -857.4000000000001
$ws.Cells.Item(22, 14).Value = -1881.9333  # N22: -1948.5 -> -1881.9333
$ws.Cells.Item(58, 8).Value = 8845.218000000001  # H58: 9110.454 -> 8845.218000000001
$ws.Cells.Item(58, 9).Value = 8702.875  # I58: 8278.588 -> 8702.875
$ws.Cells.Item(58, 10).Value = 9170.571  # J58: 11938.8 -> 9170.571
$ws.Cells.Item(58, 11).Value = 8702.875  # K58: 8278.588 -> 8702.875
$ws.Cells.Item(58, 12).Value = 9170.571  # L58: 11938.8 -> 9170.571
$ws.Cells.Item(58, 13).Value = -8499.875  # M58: -8075.588 -> -8499.875
$ws.Cells.Item(58, 14).Value = -9576.571  # N58: -12344.8 -> -9576.571
$ws.Cells.Item(134, 8).Value = 7405.564  # H134: 7250.35 -> 7405.564
$ws.Cells.Item(134, 9).Value = 8361.759  # I134: 8122.933 -> 8361.759
$ws.Cells.Item(134, 11).Value = 25085.277  # K134: 24368.799 -> 25085.277
$ws.Cells.Item(134, 13).Value = -22550.277  # M134: -21833.799 -> -22550.277
$ws.Cells.Item(136, 8).Value = 8845.218000000001  # H136: 9110.454 -> 8845.218000000001
$ws.Cells.Item(136, 9).Value = 8702.875  # I136: 8278.588 -> 8702.875
$ws.Cells.Item(136, 10).Value = 9170.571  # J136: 11938.8 -> 9170.571
$ws.Cells.Item(136, 11).Value = 26108.625  # K136: 24835.764 -> 26108.625
$ws.Cells.Item(136, 12).Value = 27511.713  # L136: 35816.39999999999 -> 27511.713
$ws.Cells.Item(136, 13).Value = -23558.625  # M136: -22285.764 -> -23558.625
$ws.Cells.Item(136, 14).Value = -32611.713  # N136: -40916.39999999999 -> -32611.713

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 84208.664  # H55: 87870.78 -> 84208.664
$ws.Cells.Item(55, 10).Value = 91748.09  # J55: 96118 -> 91748.09
$ws.Cells.Item(55, 12).Value = 275244.27  # L55: 288354 -> 275244.27
$ws.Cells.Item(55, 14).Value = -275598.27  # N55: -288708 -> -275598.27
$ws.Cells.Item(63, 8).Value = 12800  # H63: 15000 -> 12800
$ws.Cells.Item(63, 9).Value = 8400  # I63: 0 -> 8400
$ws.Cells.Item(63, 11).Value = 25200  # K63: 0 -> 25200
$ws.Cells.Item(63, 13).Value = -24451  # M63: None -> -24451
$ws.Cells.Item(66, 8).Value = 12800  # H66: 15000 -> 12800
$ws.Cells.Item(66, 9).Value = 8400  # I66: 0 -> 8400
$ws.Cells.Item(66, 11).Value = 75600  # K66: 0 -> 75600
$ws.Cells.Item(66, 13).Value = -71856  # M66: None -> -71856
$ws.Cells.Item(131, 8).Value = 1549614.4  # H131: 1472158.6 -> 1549614.4
$ws.Cells.Item(131, 9).Value = 3678086.5  # I131: 3269465.8 -> 3678086.5
$ws.Cells.Item(131, 11).Value = 11034259.5  # K131: 9808397.399999999 -> 11034259.5
$ws.Cells.Item(131, 13).Value = -11029219.5  # M131: -9803357.399999999 -> -11029219.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3555.8333  # H122: 3606.36 -> 3555.8333
$ws.Cells.Item(122, 9).Value = 3488.9285  # I122: 3577.6 -> 3488.9285
$ws.Cells.Item(122, 11).Value = 10466.7855  # K122: 10732.8 -> 10466.7855
$ws.Cells.Item(122, 13).Value = -8016.7855  # M122: -8282.799999999999 -> -8016.7855
$ws.Cells.Item(126, 8).Value = 5327.3  # H126: 5361.846 -> 5327.3
$ws.Cells.Item(126, 9).Value = 6876.25  # I126: 7290 -> 6876.25
$ws.Cells.Item(126, 11).Value = 20628.75  # K126: 21870 -> 20628.75
$ws.Cells.Item(126, 13).Value = -18158.75  # M126: -19400 -> -18158.75
$ws.Cells.Item(132, 8).Value = 4185.054  # H132: 4440.206 -> 4185.054
$ws.Cells.Item(132, 9).Value = 3847.7666  # I132: 4131.593 -> 3847.7666
$ws.Cells.Item(132, 11).Value = 11543.2998  # K132: 12394.779 -> 11543.2998
$ws.Cells.Item(132, 13).Value = -9013.299800000001  # M132: -9864.778999999999 -> -9013.299800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2026.375  # H16: 6946467.5 -> 2026.375
$ws.Cells.Item(16, 9).Value = 1894.6666  # I16: 7814276 -> 1894.6666
$ws.Cells.Item(16, 10).Value = 4002  # J16: 3997 -> 4002
$ws.Cells.Item(16, 11).Value = 1894.6666  # K16: 7814276 -> 1894.6666
$ws.Cells.Item(16, 12).Value = 4002  # L16: 3997 -> 4002
$ws.Cells.Item(16, 13).Value = -1724.6666  # M16: -7814106 -> -1724.6666
$ws.Cells.Item(16, 14).Value = -4342  # N16: -4337 -> -4342
$ws.Cells.Item(55, 8).Value = 1213.3889  # H55: 1228.0555 -> 1213.3889
$ws.Cells.Item(55, 10).Value = 1726.2858  # J55: 1764 -> 1726.2858
$ws.Cells.Item(55, 12).Value = 1726.2858  # L55: 1764 -> 1726.2858
$ws.Cells.Item(55, 14).Value = -2072.2858  # N55: -2110 -> -2072.2858
$ws.Cells.Item(93, 8).Value = 1314.5714  # H93: 1350.7693 -> 1314.5714
$ws.Cells.Item(93, 9).Value = 1314.5714  # I93: 1380.25 -> 1314.5714
$ws.Cells.Item(93, 10).Value = 0  # J93: 997 -> 0
$ws.Cells.Item(93, 11).Value = 1314.5714  # K93: 1380.25 -> 1314.5714
$ws.Cells.Item(93, 12).Value = 0  # L93: 997 -> 0
$ws.Cells.Item(93, 13).Value = -66.57140000000004  # M93: -132.25 -> -66.57140000000004
$ws.Cells.Item(93, 14).ClearContents()  # N93: remove (was -3493)
$ws.Cells.Item(122, 8).Value = 9055.1  # H122: 9436.655000000001 -> 9055.1
$ws.Cells.Item(122, 9).Value = 9045.478999999999  # I122: 9548 -> 9045.478999999999
$ws.Cells.Item(122, 11).Value = 27136.437  # K122: 28644 -> 27136.437
$ws.Cells.Item(122, 13).Value = -24686.437  # M122: -26194 -> -24686.437
$ws.Cells.Item(132, 8).Value = 10065.25  # H132: 36820.35 -> 10065.25
$ws.Cells.Item(132, 9).Value = 10065.25  # I132: 36820.35 -> 10065.25
$ws.Cells.Item(132, 11).Value = 30195.75  # K132: 110461.05 -> 30195.75
$ws.Cells.Item(132, 13).Value = -27665.75  # M132: -107931.05 -> -27665.75
$ws.Cells.Item(136, 8).Value = 5691.3335  # H136: 7051.0547 -> 5691.3335
$ws.Cells.Item(136, 9).Value = 2572.7273  # I136: 5485.5884 -> 2572.7273
$ws.Cells.Item(136, 10).Value = 12552.267  # J136: 9585.619000000001 -> 12552.267
$ws.Cells.Item(136, 11).Value = 7718.1819  # K136: 16456.7652 -> 7718.1819
$ws.Cells.Item(136, 12).Value = 37656.801  # L136: 28756.857 -> 37656.801
$ws.Cells.Item(136, 13).Value = -5168.1819  # M136: -13906.7652 -> -5168.1819
$ws.Cells.Item(136, 14).Value = -42756.801  # N136: -33856.857 -> -42756.801

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 8120.8335  # H5: 8122.5 -> 8120.8335
$ws.Cells.Item(5, 10).Value = 4945  # J5: 4947 -> 4945
$ws.Cells.Item(5, 12).Value = 4945  # L5: 4947 -> 4945
$ws.Cells.Item(5, 14).Value = -5169  # N5: -5171 -> -5169
$ws.Cells.Item(104, 8).Value = 20309.715  # H104: 17719.834 -> 20309.715
$ws.Cells.Item(104, 10).Value = 20309.715  # J104: 17719.834 -> 20309.715
$ws.Cells.Item(104, 12).Value = 20309.715  # L104: 17719.834 -> 20309.715
$ws.Cells.Item(104, 14).Value = -27297.715  # N104: -24707.834 -> -27297.715
$ws.Cells.Item(107, 8).Value = 9804248  # H107: 10417004 -> 9804248
$ws.Cells.Item(107, 9).Value = 324.66666  # I107: 346.25 -> 324.66666
$ws.Cells.Item(107, 11).Value = 973.9999799999999  # K107: 1038.75 -> 973.9999799999999
$ws.Cells.Item(107, 13).Value = 946.0000200000001  # M107: 881.25 -> 946.0000200000001
$ws.Cells.Item(132, 8).Value = 2274.94  # H132: 2239.8235 -> 2274.94
$ws.Cells.Item(132, 9).Value = 2305.6736  # I132: 2269.24 -> 2305.6736
$ws.Cells.Item(132, 11).Value = 6917.0208  # K132: 6807.719999999999 -> 6917.0208
$ws.Cells.Item(132, 13).Value = -4387.0208  # M132: -4277.719999999999 -> -4387.0208
